# Regenerate orders with updated distance/size codes.
#
# The experiment's distance and size condition labels were updated:
#   D80 -> D86, D51 -> D55, D64 -> D69, S30 -> S31
#
# These codes appear as substrings inside many cell values throughout the
# sheet (Condition, Filename_Left, Filename_Right, Distance, Size columns),
# exactly as they would if the order file had been regenerated from an
# updated parameter table. The most faithful reproduction is therefore a
# literal substring find & replace applied across every used cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPart = 2       # Excel's xlPart constant -> LookAt:=xlPart (substring match)
$xlByRows = 1      # Excel's xlByRows constant -> SearchOrder:=xlByRows

$rng = $ws.UsedRange

# These four old/new substrings are mutually exclusive (none is a substring
# of another), so the four replacements are independent and order-safe.
$rng.Replace("D80", "D86", $xlPart, $xlByRows, $false, $false, $true, $true)
$rng.Replace("D51", "D55", $xlPart, $xlByRows, $false, $false, $true, $true)
$rng.Replace("D64", "D69", $xlPart, $xlByRows, $false, $false, $true, $true)
$rng.Replace("S30", "S31", $xlPart, $xlByRows, $false, $false, $true, $true)
